$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "61.028.36"
$dCell.ClearFormats()
$ws.Range("E2").Value = "  +0.24%  "

$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "2.926.13"
$dCell.ClearFormats()
$ws.Range("E3").Value = "  +0.23%  "

$dCell = $ws.Range("D4")
$dCell.NumberFormat = "@"
$dCell.Value = "0.999"
$dCell.ClearFormats()
$ws.Range("E4").Value = "  -0.05%  "

$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "591.08"
$dCell.ClearFormats()
$ws.Range("E5").Value = "  +1.10%  "

$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "147.13"
$dCell.ClearFormats()
$ws.Range("E6").Value = "  +1.17%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  +0.56%  "

$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = "6.93"
$dCell.ClearFormats()
$ws.Range("E9").Value = "  +0.89%  "

$ws.Range("E10").Value = "  -0.42%  "

$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = "0.441"
$dCell.ClearFormats()
$ws.Range("E11").Value = "  -1.28%  "

$ws.Range("E12").Value = "  +0.04%  "

$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "33.72"
$dCell.ClearFormats()
$ws.Range("E13").Value = "  -0.01%  "

$ws.Range("E14").Value = "  -0.04%  "

$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "3.409.07"
$dCell.ClearFormats()
$ws.Range("E15").Value = "  +0.23%  "

$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "60.965.79"
$dCell.ClearFormats()
$ws.Range("E16").Value = "  +0.20%  "

$ws.Range("E17").Value = "  -0.57%  "

$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "2.923.67"
$dCell.ClearFormats()
$ws.Range("E18").Value = "  +0.23%  "

$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = "431.73"
$dCell.ClearFormats()
$ws.Range("E19").Value = "  -0.03%  "

$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "13.41"
$dCell.ClearFormats()
$ws.Range("E20").Value = "  -1.70%  "

$ws.Range("E21").Value = "  -0.59%  "

$ws.Range("E22").Value = "  -0.69%  "

$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "81.39"
$dCell.ClearFormats()
$ws.Range("E23").Value = "  +1.20%  "

$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = "10.94"
$dCell.ClearFormats()
$ws.Range("E24").Value = "  +0.61%  "

$ws.Range("E25").Value = "  -0.46%  "

$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = "11.91"
$dCell.ClearFormats()
$ws.Range("E26").Value = "  -0.21%  "

$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("E28").Value = "  +4.90%  "

$ws.Range("E29").Value = "  +0.21%  "

$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = "7.03"
$dCell.ClearFormats()
$ws.Range("E30").Value = "  -2.89%  "

$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = "26.69"
$dCell.ClearFormats()
$ws.Range("E31").Value = "  +0.44%  "

$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "0.110"
$dCell.ClearFormats()
$ws.Range("E32").Value = "  +2.44%  "

$ws.Range("E33").Value = "  +0.03%  "

$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0₃0866"
$dCell.ClearFormats()
$ws.Range("E34").Value = "  -1.09%  "

$ws.Range("E35").Value = "  +0.08%  "

$ws.Range("E36").Value = "  -0.15%  "

$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = "3.02"
$dCell.ClearFormats()
$ws.Range("E37").Value = "  -0.85%  "

$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = "1.99"
$dCell.ClearFormats()
$ws.Range("E38").Value = "  -0.99%  "

$ws.Range("E39").Value = "  -4.89%  "

$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "8.57"
$dCell.ClearFormats()
$ws.Range("E40").Value = "  -1.08%  "

$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = "41.52"
$dCell.ClearFormats()
$ws.Range("E41").Value = "  +0.99%  "

$ws.Range("E42").Value = "  -4.68%  "

$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = "378.46"
$dCell.ClearFormats()
$ws.Range("E43").Value = "  +0.07%  "

$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = "2.709.06"
$dCell.ClearFormats()
$ws.Range("E44").Value = "  +1.01%  "

$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0344"
$dCell.ClearFormats()
$ws.Range("E45").Value = "  -1.69%  "

$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = "134.02"
$dCell.ClearFormats()
$ws.Range("E46").Value = "  +1.12%  "

$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = "23.87"
$dCell.ClearFormats()
$ws.Range("E48").Value = "  -2.90%  "

$ws.Range("E49").Value = "  -0.57%  "

$ws.Range("E50").Value = "  -2.62%  "

$ws.Range("E51").Value = "  -0.65%  "
